# F360 Studyplan.xlsx — editor script
#
# Changes applied (per commit "added collabs and improved excel"):
#  1. Row 12 was an exact duplicate of row 11 (Day 9 appeared twice under
#     "Week 2"). Delete that duplicate row outright — this naturally shifts
#     every row below it up by one (so the old row 13 "Day 10" becomes the
#     new row 12, etc.), shrinks the sheet's used dimension by one row, and
#     shrinks Table1's range/autofilter by one row as well.
#  2. The exercise example text in D7 ("Week 1" / Day 5, Patterns & Arrays)
#     is simplified from mentioning two projects to just one.
#  3. The active selection is moved to the newly-adjacent row (the whole of
#     row 11) to match where the author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the duplicated "Day 9" row — shifts rows 13:31 up to 12:30,
#    and Excel auto-adjusts dimension / Table1 ref / autoFilter accordingly.
$ws.Rows(12).Delete()

# 2) Update the exercise example text that no longer mentions the phone stand.
$ws.Range("D7").Value = "Design a gear wheel with patterned holes"

# 3) Reflect the author's final on-screen selection (whole row 11).
$ws.Range("A11:XFD11").Select() | Out-Null
